$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 118.0346986666667
$ws.Range("H2").Value = 354.104096
$ws.Range("I2").Value = 0.2666057129183408
$ws.Range("J2").Value = 0.2666057129183408
$ws.Range("M2").Value = 9.467112666666667
$ws.Range("N2").Value = 28.401338
$ws.Range("O2").Value = 0.2314516669582087
$ws.Range("P2").Value = 0.2314516669582087
$ws.Range("Q2").Value = 1117.447790853383
$ws.Range("R2").Value = 10057.03011768045
$ws.Range("S2").Value = 0.06170633667553162
$ws.Range("T2").Value = 0.06170633667553162
$ws.Range("G3").Value = 118.0346986666667
$ws.Range("H3").Value = 354.104096
$ws.Range("I3").Value = 0.2666057129183408
$ws.Range("J3").Value = 0.2666057129183408
$ws.Range("O3").Value = 0.03494502079849753
$ws.Range("P3").Value = 0.03494502079849753
$ws.Range("Q3").Value = 168.7144309902827
$ws.Range("R3").Value = 1518.429878912544
$ws.Range("S3").Value = 0.009316542182929681
$ws.Range("T3").Value = 0.009316542182929681
$ws.Range("G4").Value = 118.0346986666667
$ws.Range("H4").Value = 354.104096
$ws.Range("I4").Value = 0.2666057129183408
$ws.Range("J4").Value = 0.2666057129183408
$ws.Range("M4").Value = 10.196198
$ws.Range("N4").Value = 30.588594
$ws.Range("O4").Value = 0.2492763218130026
$ws.Range("P4").Value = 0.2492763218130026
$ws.Range("Q4").Value = 1203.505158475669
$ws.Range("R4").Value = 10831.54642628102
$ws.Range("S4").Value = 0.06645849149061732
$ws.Range("T4").Value = 0.06645849149061732
$ws.Range("G5").Value = 118.0346986666667
$ws.Range("H5").Value = 354.104096
$ws.Range("I5").Value = 0.2666057129183408
$ws.Range("J5").Value = 0.2666057129183408
$ws.Range("M5").Value = 19.81052133333333
$ws.Range("N5").Value = 59.431564
$ws.Range("O5").Value = 0.4843269904302911
$ws.Range("P5").Value = 0.4843269904302911
$ws.Range("Q5").Value = 2338.328916009571
$ws.Range("R5").Value = 21044.96024408614
$ws.Range("S5").Value = 0.1291243425692622
$ws.Range("T5").Value = 0.1291243425692622
$ws.Range("I6").Value = 0.4881754016778185
$ws.Range("J6").Value = 0.4881754016778186
$ws.Range("M6").Value = 9.467112666666667
$ws.Range("N6").Value = 28.401338
$ws.Range("O6").Value = 0.2314516669582087
$ws.Range("P6").Value = 0.2314516669582087
$ws.Range("Q6").Value = 2046.132163420394
$ws.Range("R6").Value = 18415.18947078355
$ws.Range("S6").Value = 0.1129890104863242
$ws.Range("T6").Value = 0.1129890104863242
$ws.Range("I7").Value = 0.4881754016778185
$ws.Range("J7").Value = 0.4881754016778186
$ws.Range("O7").Value = 0.03494502079849753
$ws.Range("P7").Value = 0.03494502079849753
$ws.Range("S7").Value = 0.01705929956494625
$ws.Range("T7").Value = 0.01705929956494625
$ws.Range("I8").Value = 0.4881754016778185
$ws.Range("J8").Value = 0.4881754016778186
$ws.Range("M8").Value = 10.196198
$ws.Range("N8").Value = 30.588594
$ws.Range("O8").Value = 0.2492763218130026
$ws.Range("P8").Value = 0.2492763218130026
$ws.Range("Q8").Value = 2203.709769490722
$ws.Range("R8").Value = 19833.3879254165
$ws.Range("S8").Value = 0.1216905685298317
$ws.Range("T8").Value = 0.1216905685298317
$ws.Range("I9").Value = 0.4881754016778185
$ws.Range("J9").Value = 0.4881754016778186
$ws.Range("M9").Value = 19.81052133333333
$ws.Range("N9").Value = 59.431564
$ws.Range("O9").Value = 0.4843269904302911
$ws.Range("P9").Value = 0.4843269904302911
$ws.Range("Q9").Value = 4281.658653644332
$ws.Range("R9").Value = 38534.92788279899
$ws.Range("S9").Value = 0.2364365230967164
$ws.Range("T9").Value = 0.2364365230967164
$ws.Range("G10").Value = 45.876452
$ws.Range("H10").Value = 137.629356
$ws.Range("I10").Value = 0.1036214293744632
$ws.Range("J10").Value = 0.1036214293744632
$ws.Range("M10").Value = 9.467112666666667
$ws.Range("N10").Value = 28.401338
$ws.Range("O10").Value = 0.2314516669582087
$ws.Range("P10").Value = 0.2314516669582087
$ws.Range("Q10").Value = 434.3175398309253
$ws.Range("R10").Value = 3908.857858478329
$ws.Range("S10").Value = 0.0239833525613118
$ws.Range("T10").Value = 0.0239833525613118
$ws.Range("G11").Value = 45.876452
$ws.Range("H11").Value = 137.629356
$ws.Range("I11").Value = 0.1036214293744632
$ws.Range("J11").Value = 0.1036214293744632
$ws.Range("O11").Value = 0.03494502079849753
$ws.Range("P11").Value = 0.03494502079849753
$ws.Range("Q11").Value = 65.57410306007601
$ws.Range("R11").Value = 590.166927540684
$ws.Range("S11").Value = 0.003621053004660659
$ws.Range("T11").Value = 0.003621053004660659
$ws.Range("G12").Value = 45.876452
$ws.Range("H12").Value = 137.629356
$ws.Range("I12").Value = 0.1036214293744632
$ws.Range("J12").Value = 0.1036214293744632
$ws.Range("M12").Value = 10.196198
$ws.Range("N12").Value = 30.588594
$ws.Range("O12").Value = 0.2492763218130026
$ws.Range("P12").Value = 0.2492763218130026
$ws.Range("Q12").Value = 467.765388129496
$ws.Range("R12").Value = 4209.888493165464
$ws.Range("S12").Value = 0.02583036877547201
$ws.Range("T12").Value = 0.02583036877547201
$ws.Range("G13").Value = 45.876452
$ws.Range("H13").Value = 137.629356
$ws.Range("I13").Value = 0.1036214293744632
$ws.Range("J13").Value = 0.1036214293744632
$ws.Range("M13").Value = 19.81052133333333
$ws.Range("N13").Value = 59.431564
$ws.Range("O13").Value = 0.4843269904302911
$ws.Range("P13").Value = 0.4843269904302911
$ws.Range("Q13").Value = 908.8364310436427
$ws.Range("R13").Value = 8179.527879392785
$ws.Range("S13").Value = 0.05018665503301872
$ws.Range("T13").Value = 0.05018665503301873
$ws.Range("G14").Value = 62.68962833333333
$ws.Range("H14").Value = 188.068885
$ws.Range("I14").Value = 0.1415974560293775
$ws.Range("J14").Value = 0.1415974560293775
$ws.Range("M14").Value = 9.467112666666667
$ws.Range("N14").Value = 28.401338
$ws.Range("O14").Value = 0.2314516669582087
$ws.Range("P14").Value = 0.2314516669582087
$ws.Range("Q14").Value = 593.4897744631255
$ws.Range("R14").Value = 5341.40797016813
$ws.Range("S14").Value = 0.03277296723504107
$ws.Range("T14").Value = 0.03277296723504108
$ws.Range("G15").Value = 62.68962833333333
$ws.Range("H15").Value = 188.068885
$ws.Range("I15").Value = 0.1415974560293775
$ws.Range("J15").Value = 0.1415974560293775
$ws.Range("O15").Value = 0.03494502079849753
$ws.Range("P15").Value = 0.03494502079849753
$ws.Range("Q15").Value = 89.60623522341834
$ws.Range("R15").Value = 806.4561170107651
$ws.Range("S15").Value = 0.004948126045960934
$ws.Range("T15").Value = 0.004948126045960935
$ws.Range("G16").Value = 62.68962833333333
$ws.Range("H16").Value = 188.068885
$ws.Range("I16").Value = 0.1415974560293775
$ws.Range("J16").Value = 0.1415974560293775
$ws.Range("M16").Value = 10.196198
$ws.Range("N16").Value = 30.588594
$ws.Range("O16").Value = 0.2492763218130026
$ws.Range("P16").Value = 0.2492763218130026
$ws.Range("Q16").Value = 639.1958630330766
$ws.Range("R16").Value = 5752.76276729769
$ws.Range("S16").Value = 0.03529689301708158
$ws.Range("T16").Value = 0.03529689301708159
$ws.Range("G17").Value = 62.68962833333333
$ws.Range("H17").Value = 188.068885
$ws.Range("I17").Value = 0.1415974560293775
$ws.Range("J17").Value = 0.1415974560293775
$ws.Range("M17").Value = 19.81052133333333
$ws.Range("N17").Value = 59.431564
$ws.Range("O17").Value = 0.4843269904302911
$ws.Range("P17").Value = 0.4843269904302911
$ws.Range("Q17").Value = 1241.914219476238
$ws.Range("R17").Value = 11177.22797528614
$ws.Range("S17").Value = 0.06857946973129386
$ws.Range("T17").Value = 0.06857946973129388
